# Sexta versión — blank out the personal-data placeholders in the tutela
# template (judge, accionante/accionado names, cédula number, city,
# address and petition date) with underscore blanks, per the commit.
#
# Each entry below is the *exact* text of one (or more identical) source
# run(s), longest/most-specific first. Longer, more context-qualified
# strings are replaced before the short bare strings ("Radamel Falcao",
# "Nicolás", "1020829275") so a short replacement can't fire early and
# clobber a longer match that still contains it.

$d = $word.ActiveDocument

$replacements = @(
    @{Old = "En el presente caso, esta legitimación se encuentra demostrada toda vez que soy el titular del derecho fundamental de petición, que, en este caso, fue vulnerado por Radamel Falcao";
      New = "En el presente caso, esta legitimación se encuentra demostrada toda vez que soy el titular del derecho fundamental de petición, que, en este caso, fue vulnerado por ____________"},
    @{Old = " dar una respuesta clara, de fondo y sin respuestas evasivas al derecho de petición radicado el día 2022-07-11";
      New = " dar una respuesta clara, de fondo y sin respuestas evasivas al derecho de petición radicado el día _____________"},
    @{Old = " radiqué un derecho de petición de solicitud de documentos ante Radamel Falcao";
      New = " radiqué un derecho de petición de solicitud de documentos ante ____________"},
    @{Old = "tición (artículo 23 Constitución Política), en contra de: Radamel Falcao";
      New = "tición (artículo 23 Constitución Política), en contra de: ____________"},
    @{Old = "En el presente caso, fue la autoridad pública Radamel Falcao";
      New = "En el presente caso, fue la autoridad pública ____________"},
    @{Old = ", identificado con cédula de ciudadanía número 1020829275";
      New = ", identificado con cédula de ciudadanía número ______________"},
    @{Old = "1. Que se declare que con el actuar de Radamel Falcao";
      New = "1. Que se declare que con el actuar de ____________"},
    @{Old = "2. Derecho de petición radicado ante Radamel Falcao";
      New = "2. Derecho de petición radicado ante ____________"},
    @{Old = ", o quien haga sus veces, con domicilio en fusa";
      New = ", o quien haga sus veces, con domicilio en _____________"},
    @{Old = ", domiciliado en la ciudad de Bogotá";
      New = ", domiciliado en la ciudad de _____________"},
    @{Old = "2. Que se ordene a Radamel Falcao";
      New = "2. Que se ordene a ____________"},
    @{Old = "3. Que se ordene a Radamel Falcao";
      New = "3. Que se ordene a ____________"},
    @{Old = " el día 2022-07-11";
      New = " el día _____________"},
    @{Old = "1. El 2022-07-11";
      New = "1. El _____________"},
    @{Old = "Radamel Falcao";
      New = "____________"},
    @{Old = "JUEZ(A) de ";
      New = "JUEZ(A) de __________"},
    @{Old = "1020829275";
      New = "______________"},
    @{Old = "Nicolás";
      New = "____________"}
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $pair.New, 2) | Out-Null
}

Write-Output "Applied $($replacements.Count) placeholder replacements."
